# CryCompanywiseStockReport_1.xlsx -- stock-count correction pass.
#
# Several duplicate-SKU row pairs had their Closing Qty/Value (and the
# derived rate columns) entered against the wrong one of the two rows --
# this swaps B (item code), E (rate), F (qty) and G (value) between each
# such pair. A separate batch of rows simply had their Closing Qty (F)
# reduced by a small count, with Closing Value (G = Rate(D) x Qty(F))
# recalculated to match. Company "Sub Total:" rows and the workbook
# "Sub Total:"/"Grand Total:" rows are then refreshed to the resulting
# sum of Closing Value for their block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $value) {
    $ws.Range($addr).Value2 = $value
}

# --- Row pairs: swap Item Code (B), Rate (E), Qty (F) and Value (G) ----
# Rows 127 / 128
Set-Cell "B127" 64329
Set-Cell "E127" 128.32
Set-Cell "F127" 2
Set-Cell "G127" 241.38
Set-Cell "B128" 57552
Set-Cell "E128" 136.86
Set-Cell "F128" -5
Set-Cell "G128" -603.45

# Rows 192 / 193
Set-Cell "B192" 48706
Set-Cell "E192" 39.8
Set-Cell "F192" -144
Set-Cell "G192" -4795.2
Set-Cell "B193" 64973
Set-Cell "E193" 35.4
Set-Cell "F193" 2
Set-Cell "G193" 66.59999999999999

# Rows 243 / 244
Set-Cell "B243" 60325
Set-Cell "E243" 151.57
Set-Cell "F243" -102
Set-Cell "G243" -12939.72
Set-Cell "B244" 63560
Set-Cell "E244" 134.87
Set-Cell "F244" 1
Set-Cell "G244" 126.86

# Rows 366 / 367
Set-Cell "B366" 53263
Set-Cell "E366" 15.29
Set-Cell "F366" -309
Set-Cell "G366" -3958.29
Set-Cell "B367" 65066
Set-Cell "E367" 13.61
Set-Cell "F367" 90
Set-Cell "G367" 1152.9

# Rows 375 / 376
Set-Cell "B375" 64927
Set-Cell "E375" 17.26
Set-Cell "F375" 106
Set-Cell "G375" 1719.32
Set-Cell "B376" 45718
Set-Cell "E376" 19.38
Set-Cell "F376" -294
Set-Cell "G376" -4768.68

# Rows 380 / 381
Set-Cell "B380" 45709
Set-Cell "E380" 15.69
Set-Cell "F380" -300
Set-Cell "G380" -3945
Set-Cell "B381" 64925
Set-Cell "E381" 13.97
Set-Cell "F381" 111
Set-Cell "G381" 1459.65

# Rows 442 / 443
Set-Cell "B442" 53319
Set-Cell "E442" 310.64
Set-Cell "F442" -6
Set-Cell "G442" -1643.52
Set-Cell "B443" 64810
Set-Cell "E443" 291.22
Set-Cell "F443" 5
Set-Cell "G443" 1369.6

# Rows 572 / 573
Set-Cell "B572" 65079
Set-Cell "F572" 18
Set-Cell "G572" 735.66
Set-Cell "B573" 65362
Set-Cell "F573" 30
Set-Cell "G573" 1226.1

# --- Single rows: reduce Closing Qty (F), recompute Closing Value (G) ----
Set-Cell "F9" 18
Set-Cell "G9" 532.26
Set-Cell "F77" 309
Set-Cell "G77" 14442.66
Set-Cell "F102" 8
Set-Cell "G102" 395.84
Set-Cell "F115" 240
Set-Cell "G115" 23234.4
Set-Cell "F149" 258
Set-Cell "G149" 16718.4
Set-Cell "F150" 50
Set-Cell "G150" 2324.5
Set-Cell "F167" 23
Set-Cell "G167" 6601.23
Set-Cell "F186" 82
Set-Cell "G186" 3523.54
Set-Cell "F229" 72
Set-Cell "G229" 10330.56
Set-Cell "F234" 47
Set-Cell "G234" 2412.04
Set-Cell "F249" 148
Set-Cell "G249" 20397.36
Set-Cell "F255" 629
Set-Cell "G255" 107766.57
Set-Cell "F256" 305
Set-Cell "G256" 46106.85
Set-Cell "F273" 21
Set-Cell "G273" 6741.21
Set-Cell "F278" 19
Set-Cell "G278" 2609.08
Set-Cell "F280" 150
Set-Cell "G280" 25371
Set-Cell "F282" 18
Set-Cell "G282" 966.6
Set-Cell "F285" 35
Set-Cell "G285" 977.55
Set-Cell "F294" 55
Set-Cell "G294" 3924.8
Set-Cell "F295" 8
Set-Cell "G295" 829.52
Set-Cell "F296" 107
Set-Cell "G296" 2268.4
Set-Cell "F301" 10
Set-Cell "G301" 5061.4
Set-Cell "F307" 2
Set-Cell "G307" 308.3
Set-Cell "F336" 23
Set-Cell "G336" 1003.95
Set-Cell "F338" 87
Set-Cell "G338" 2061.9
Set-Cell "F341" 7
Set-Cell "G341" 356.65
Set-Cell "F345" 95
Set-Cell "G345" 5833.95
Set-Cell "F353" 22
Set-Cell "G353" 3018.18
Set-Cell "F354" 28
Set-Cell "G354" 1920.52
Set-Cell "F390" 12
Set-Cell "G390" 739.8
Set-Cell "F402" 1
Set-Cell "G402" 16.18
Set-Cell "F429" 24
Set-Cell "G429" 163.2
Set-Cell "F430" 18
Set-Cell "G430" 232.02
Set-Cell "F453" 35
Set-Cell "G453" 927.85
Set-Cell "F490" 11
Set-Cell "G490" 1436.05
Set-Cell "F509" 265
Set-Cell "G509" 21300.7
Set-Cell "F542" 55
Set-Cell "G542" 7124.15
Set-Cell "F551" 30
Set-Cell "G551" 4293.9
Set-Cell "F552" 36
Set-Cell "G552" 3664.44
Set-Cell "F553" 0
Set-Cell "G553" 0
Set-Cell "F555" 44
Set-Cell "G555" 3060.64
Set-Cell "F558" 80
Set-Cell "G558" 10800.8
Set-Cell "F577" 88
Set-Cell "G577" 3783.12
Set-Cell "F579" 39
Set-Cell "G579" 3143.4
Set-Cell "F580" 76
Set-Cell "G580" 4331.24
Set-Cell "F581" 35
Set-Cell "G581" 8463
Set-Cell "F582" 65
Set-Cell "G582" 3704.35
Set-Cell "F599" 2355
Set-Cell "G599" 384124.05
Set-Cell "F601" 499
Set-Cell "G601" 141152.13
Set-Cell "F602" 380
Set-Cell "G602" 54967
Set-Cell "F613" 162
Set-Cell "G613" 25783.92
Set-Cell "F617" 16
Set-Cell "G617" 632.48

# --- Sub Total: / Grand Total: rows: refresh rolled-up Closing Value (B) ----
Set-Cell "B10" 33757.61
Set-Cell "B90" 212516.53
Set-Cell "B104" 363.82
Set-Cell "B117" 17401.86
Set-Cell "B156" 37382.23
Set-Cell "B175" 37592.11
Set-Cell "B216" 57253.22
Set-Cell "B260" 222377.16
Set-Cell "B275" 9756.93
Set-Cell "B304" 206074.98
Set-Cell "B309" 2397.76
Set-Cell "B346" 31546.39
Set-Cell "B358" 39218.08
Set-Cell "B395" 1756.91
Set-Cell "B411" 10563.06
Set-Cell "B435" 2245.12
Set-Cell "B460" 17078.55
Set-Cell "B493" 15820.21
Set-Cell "B510" 28784.38
Set-Cell "B547" 25780.59
Set-Cell "B560" 25143.78
Set-Cell "B583" 34407.27
Set-Cell "B606" 581091.23
Set-Cell "B618" 51232.72
Set-Cell "B619" 2157496.59
Set-Cell "B620" 2157496.59
